$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hardware")

# --- Fill in new BOM rows 8-27 (Motor Control PCB parts) ---
# Row 8
$ws.Range("A8").Value = "PFET"
$ws.Range("B8").Value = "FQP27P06"
$ws.Range("B8").HorizontalAlignment = -4131
$ws.Range("B8").VerticalAlignment = -4108
$ws.Range("C8").Value = "Mouser"
$ws.Range("D8").Value = "Motor Control PCB"
$ws.Range("E8").Value = "http://www.mouser.com/ds/2/149/FQP27P06-244293.pdf"
$ws.Range("F8").Value = 1.28
$ws.Range("G8").Value = 20
$ws.Range("H8").Formula = "=F8*G8"

# Row 9
$ws.Range("A9").Value = "NFET"
$ws.Range("B9").Value = "FQP30N06L"
$ws.Range("B9").HorizontalAlignment = -4131
$ws.Range("B9").VerticalAlignment = -4108
$ws.Range("C9").Value = "Mouser"
$ws.Range("D9").Value = "Motor Control PCB"
$ws.Range("E9").Value = "http://www.mouser.com/ds/2/149/FQP30N06L-244344.pdf"
$ws.Range("F9").Value = 1.11
$ws.Range("G9").Value = 20
$ws.Range("H9").Formula = "=F9*G9"

# Row 10
$ws.Range("A10").Value = "NPN Transisor"
$ws.Range("B10").Value = "BC846-NPN-SOT23-BEC"
$ws.Range("C10").Value = "Mouser"
$ws.Range("D10").Value = "Motor Control PCB"
$ws.Range("E10").Value = "http://www.mouser.com/ds/2/302/BC846_SER-840048.pdf"
$ws.Range("F10").Value = 0.13
$ws.Range("G10").Value = 60
$ws.Range("H10").Formula = "=F10*G10"

# Row 11
$ws.Range("A11").Value = "5V Regulator"
$ws.Range("B11").Value = "LD1085V50"
$ws.Range("C11").Value = "Mouser"
$ws.Range("D11").Value = "Motor Control PCB"
$ws.Range("E11").Value = "http://www.mouser.com/ds/2/389/CD00001883-310744.pdf"
$ws.Range("F11").Value = 1.49
$ws.Range("G11").Value = 8
$ws.Range("H11").Formula = "=F11*G11"

# Row 12
$ws.Range("A12").Value = "PTC for motors"
$ws.Range("B12").Value = "MF-SMDF050-2"
$ws.Range("C12").Value = "Mouser"
$ws.Range("D12").Value = "Motor Control PCB"
$ws.Range("E12").Value = "http://www.mouser.com/ds/2/54/fsmdf2-777754.pdf"
$ws.Range("F12").Value = 0.42
$ws.Range("G12").Value = 16
$ws.Range("H12").Formula = "=F12*G12"

# Row 13
$ws.Range("A13").Value = "Diode-SMA"
$ws.Range("B13").Value = "MURS360-E3/57T"
$ws.Range("C13").Value = "Mouser"
$ws.Range("D13").Value = "Motor Control PCB"
$ws.Range("E13").Value = "http://www.mouser.com/ds/2/427/murs340-240280.pdf"
$ws.Range("F13").Value = 0.7
$ws.Range("G13").Value = 10
$ws.Range("H13").Formula = "=F13*G13"

# Row 14
$ws.Range("A14").Value = "Diode-1N4148"
$ws.Range("B14").Value = "1N4148"
$ws.Range("C14").Value = "Mouser"
$ws.Range("D14").Value = "Motor Control PCB"
$ws.Range("E14").Value = "http://www.mouser.com/ds/2/149/1N4148-888354.pdf"
$ws.Range("F14").Value = 0.1
$ws.Range("G14").Value = 10
$ws.Range("H14").Formula = "=F14*G14"

# Row 15
$ws.Range("A15").Value = "SPOX-3 for AX12A"
$ws.Range("B15").Value = "22-43-6030"
$ws.Range("C15").Value = "Mouser"
$ws.Range("D15").Value = "Motor Control PCB"
$ws.Range("E15").Value = "http://www.mouser.com/ds/2/276/0022436030_PCB_HEADERS-227784.pdf"
$ws.Range("F15").Value = 0.47
$ws.Range("G15").Value = 15
$ws.Range("H15").Formula = "=F15*G15"

# Row 16
$ws.Range("A16").Value = "100uF Capacitor"
$ws.Range("B16").Value = "UWT1V101MCL1GS"
$ws.Range("C16").Value = "Mouser"
$ws.Range("D16").Value = "Motor Control PCB"
$ws.Range("E16").Value = "http://www.mouser.com/ds/2/293/e-uwt-880115.pdf"
$ws.Range("F16").Value = 0.36
$ws.Range("G16").Value = 15
$ws.Range("H16").Formula = "=F16*G16"

# Row 17
$ws.Range("A17").Value = "10uF Capacitor"
$ws.Range("B17").Value = "LMK316B7106KL-TD"
$ws.Range("C17").Value = "Mouser"
$ws.Range("D17").Value = "Motor Control PCB"
$ws.Range("E17").Value = "http://www.mouser.com/ds/2/396/mlcc_all_e-541559.pdf"
$ws.Range("F17").Value = 0.29
$ws.Range("G17").Value = 15
$ws.Range("H17").Formula = "=F17*G17"

# Row 18
$ws.Range("A18").Value = "0.1uF Capacitor"
$ws.Range("B18").Value = "GRM155R61A104KA01D"
$ws.Range("C18").Value = "Mouser"
$ws.Range("D18").Value = "Motor Control PCB"
$ws.Range("E18").Value = "http://www.mouser.com/ds/2/281/c02e-2905.pdf"
$ws.Range("F18").Value = 0.1
$ws.Range("G18").Value = 15
$ws.Range("H18").Formula = "=F18*G18"

# Row 19
$ws.Range("A19").Value = "1.7K Resistor"
$ws.Range("B19").Value = "TNPW08051K69BEEA"
$ws.Range("C19").Value = "Mouser"
$ws.Range("D19").Value = "Motor Control PCB"
$ws.Range("E19").Value = "http://www.mouser.com/ds/2/427/tnpw_e3-64594.pdf"
$ws.Range("F19").Value = 0.49
$ws.Range("G19").Value = 30
$ws.Range("H19").Formula = "=F19*G19"

# Row 20
$ws.Range("A20").Value = "3.3K Resistor"
$ws.Range("B20").Value = "TNPW08053K32BEEA"
$ws.Range("C20").Value = "Mouser"
$ws.Range("D20").Value = "Motor Control PCB"
$ws.Range("E20").Value = "http://www.mouser.com/ds/2/427/tnpw_e3-64594.pdf"
$ws.Range("F20").Value = 0.48
$ws.Range("G20").Value = 30
$ws.Range("H20").Formula = "=F20*G20"

# Row 21
$ws.Range("A21").Value = "7K Resistor"
$ws.Range("B21").Value = "RR1220P-6981-D-M"
$ws.Range("C21").Value = "Mouser"
$ws.Range("D21").Value = "Motor Control PCB"
$ws.Range("E21").Value = "http://www.mouser.com/ds/2/392/susumu_RR_Data_Sheet-358748.pdf"
$ws.Range("F21").Value = 0.1
$ws.Range("G21").Value = 10
$ws.Range("H21").Formula = "=F21*G21"

# Row 22
$ws.Range("A22").Value = "10K Resistor"
$ws.Range("B22").Value = "RR1220P-103-D"
$ws.Range("C22").Value = "Mouser"
$ws.Range("D22").Value = "Motor Control PCB"
$ws.Range("E22").Value = "http://www.mouser.com/ds/2/392/susumu_RR_Data_Sheet-358748.pdf"
$ws.Range("F22").Value = 0.1
$ws.Range("G22").Value = 25
$ws.Range("H22").Formula = "=F22*G22"

# Row 23
$ws.Range("A23").Value = "50 Resistor"
$ws.Range("B23").Value = "PCF0805-02-49R9DT1"
$ws.Range("C23").Value = "Mouser"
$ws.Range("D23").Value = "Motor Control PCB"
$ws.Range("E23").Value = "http://www.mouser.com/ds/2/414/PCF-437472.pdf"
$ws.Range("F23").Value = 0.1
$ws.Range("G23").Value = 25
$ws.Range("H23").Formula = "=F23*G23"

# Row 24
$ws.Range("A24").Value = "Molex 4 Locking"
$ws.Range("B24").Value = "22-23-2041"
$ws.Range("C24").Value = "Mouser"
$ws.Range("D24").Value = "Motor Control PCB"
$ws.Range("E24").Value = "http://www.mouser.com/ds/2/276/0022232041_PCB_HEADERS-526511.pdf"
$ws.Range("F24").Value = 0.23
$ws.Range("G24").Value = 10
$ws.Range("H24").Formula = "=F24*G24"

# Row 25
$ws.Range("A25").Value = "Molex 4 Locking 90deg"
$ws.Range("B25").Value = "22-05-3041"
$ws.Range("C25").Value = "Mouser"
$ws.Range("D25").Value = "Motor Control PCB"
$ws.Range("E25").Value = "http://www.mouser.com/ds/2/276/0022053041_PCB_HEADERS-158631.pdf"
$ws.Range("F25").Value = 0.71
$ws.Range("G25").Value = 10
$ws.Range("H25").Formula = "=F25*G25"

# Row 26
$ws.Range("A26").Value = "Molex 2 Wire Jumper assembly"
$ws.Range("A26").HorizontalAlignment = -4131
$ws.Range("A26").VerticalAlignment = -4108
$ws.Range("B26").Value = "PRT-09918"
$ws.Range("C26").Value = "Sparkfun"
$ws.Range("D26").Value = "Motor Control PCB"
$ws.Range("E26").Value = "http://www.sparkfun.com/datasheets/Prototyping/2pin_molex_set_19iv10.pdf"
$ws.Range("F26").Value = 0.95
$ws.Range("G26").Value = 10
$ws.Range("H26").Formula = "=F26*G26"

# Row 27
$ws.Range("A27").Value = "BLUE LED"
$ws.Range("A27").HorizontalAlignment = -4131
$ws.Range("A27").VerticalAlignment = -4108
$ws.Range("B27").Value = "PRT-12620"
$ws.Range("C27").Value = "Sparkfun"
$ws.Range("D27").Value = "Motor Control PCB"
$ws.Range("E27").Value = "http://cdn.sparkfun.com/datasheets/Components/LED/S150ANB4.pdf"
$ws.Range("F27").Value = 4.95
$ws.Range("G27").Value = 1
$ws.Range("I27").Value = "25 pk"
$ws.Range("H27").Formula = "=F27*G27"

# --- Hyperlinks for E column datasheet links (order matters: matches rId7..rId25) ---
$ws.Hyperlinks.Add($ws.Range("E8"), "http://www.mouser.com/ds/2/149/FQP27P06-244293.pdf")
$ws.Range("E8").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E10"), "http://www.mouser.com/ds/2/302/BC846_SER-840048.pdf")
$ws.Range("E10").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E11"), "http://www.mouser.com/ds/2/389/CD00001883-310744.pdf")
$ws.Range("E11").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E12"), "http://www.mouser.com/ds/2/54/fsmdf2-777754.pdf")
$ws.Range("E12").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E13"), "http://www.mouser.com/ds/2/427/murs340-240280.pdf")
$ws.Range("E13").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E14"), "http://www.mouser.com/ds/2/149/1N4148-888354.pdf")
$ws.Range("E14").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E15"), "http://www.mouser.com/ds/2/276/0022436030_PCB_HEADERS-227784.pdf")
$ws.Range("E15").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E16"), "http://www.mouser.com/ds/2/293/e-uwt-880115.pdf")
$ws.Range("E16").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E17"), "http://www.mouser.com/ds/2/396/mlcc_all_e-541559.pdf")
$ws.Range("E17").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E18"), "http://www.mouser.com/ds/2/281/c02e-2905.pdf")
$ws.Range("E18").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E20"), "http://www.mouser.com/ds/2/427/tnpw_e3-64594.pdf")
$ws.Range("E20").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E19"), "http://www.mouser.com/ds/2/427/tnpw_e3-64594.pdf")
$ws.Range("E19").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E21"), "http://www.mouser.com/ds/2/392/susumu_RR_Data_Sheet-358748.pdf")
$ws.Range("E21").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E22"), "http://www.mouser.com/ds/2/392/susumu_RR_Data_Sheet-358748.pdf")
$ws.Range("E22").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E23"), "http://www.mouser.com/ds/2/414/PCF-437472.pdf")
$ws.Range("E23").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E26"), "http://www.sparkfun.com/datasheets/Prototyping/2pin_molex_set_19iv10.pdf")
$ws.Range("E26").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E27"), "http://cdn.sparkfun.com/datasheets/Components/LED/S150ANB4.pdf")
$ws.Range("E27").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E25"), "http://www.mouser.com/ds/2/276/0022053041_PCB_HEADERS-158631.pdf")
$ws.Range("E25").Style = $ws.Range("E2").Style
$ws.Hyperlinks.Add($ws.Range("E24"), "http://www.mouser.com/ds/2/276/0022232041_PCB_HEADERS-526511.pdf")
$ws.Range("E24").Style = $ws.Range("E2").Style

# --- Selection matches final author cursor position ---
$ws.Range("I28").Select()
